$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DATA")

# Update K2:K5 values from 100 to 1
$ws.Range("K2:K5").Value = 1

# Update the selection to K2:K5
$ws.Range("K2:K5").Select()
